# Auto-generated edit script applying scheduled market-price refresh values
# to the Leve profit tables (Table_<Job>) across all job worksheets.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets("ALC")
$ws.Range("H40").Value = 1284.6
$ws.Range("I40").Value = 750
$ws.Range("J40").Value = 1641
$ws.Range("K40").Value = 750
$ws.Range("L40").Value = 1641
$ws.Range("M40").Value = -575
$ws.Range("N40").Value = -1991
$ws.Range("H64").Value = 3175.25
$ws.Range("I64").Value = 2002
$ws.Range("J64").Value = 3342.8572
$ws.Range("K64").Value = 2002
$ws.Range("L64").Value = 3342.8572
$ws.Range("M64").Value = -1754
$ws.Range("N64").Value = -3838.8572
$ws.Range("H67").Value = 3175.25
$ws.Range("I67").Value = 2002
$ws.Range("J67").Value = 3342.8572
$ws.Range("K67").Value = 2002
$ws.Range("L67").Value = 3342.8572
$ws.Range("M67").Value = -1144
$ws.Range("N67").Value = -5058.8572
$ws.Range("H112").Value = 1054.54
$ws.Range("J112").Value = 1065.1459
$ws.Range("L112").Value = 3195.4377
$ws.Range("N112").Value = -5411.4377
$ws.Range("H127").Value = 1005.8571
$ws.Range("I127").Value = 440.85715
$ws.Range("J127").Value = 1570.8572
$ws.Range("K127").Value = 1322.57145
$ws.Range("L127").Value = 4712.571599999999
$ws.Range("M127").Value = 3637.42855
$ws.Range("N127").Value = -14632.5716
$ws.Range("H129").Value = 176358
$ws.Range("J129").Value = 209381.77
$ws.Range("L129").Value = 628145.3099999999
$ws.Range("N129").Value = -638145.3099999999

# ---- ARM sheet ----
$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 7176.273
$ws.Range("I32").Value = 5890
$ws.Range("J32").Value = 16501.75
$ws.Range("K32").Value = 5890
$ws.Range("L32").Value = 16501.75
$ws.Range("M32").Value = -5603
$ws.Range("N32").Value = -17075.75
$ws.Range("H35").Value = 3999.6667
$ws.Range("I35").Value = 999.5
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 999.5
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = -593.5
$ws.Range("N35").Value = -10812
$ws.Range("H74").Value = 37039050
$ws.Range("J74").Value = 3415.3076
$ws.Range("L74").Value = 3415.3076
$ws.Range("N74").Value = -5163.3076
$ws.Range("H77").Value = 37039050
$ws.Range("J77").Value = 3415.3076
$ws.Range("L77").Value = 17076.538
$ws.Range("N77").Value = -25812.538
$ws.Range("H121").Value = 29857
$ws.Range("J121").Value = 29857
$ws.Range("L121").Value = 29857
$ws.Range("N121").Value = -33351
$ws.Range("H132").Value = 12077.306
$ws.Range("I132").Value = 1628.3334
$ws.Range("J132").Value = 74771.14
$ws.Range("K132").Value = 4885.0002
$ws.Range("L132").Value = 224313.42
$ws.Range("M132").Value = -2355.0002
$ws.Range("N132").Value = -229373.42
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# ---- CRP sheet ----
$ws = $wb.Worksheets("CRP")
$ws.Range("I58").Value = 1859.6364
$ws.Range("J58").Value = 73104
$ws.Range("K58").Value = 1859.6364
$ws.Range("L58").Value = 73104
$ws.Range("M58").Value = -1656.6364
$ws.Range("N58").Value = -73510
$ws.Range("H99").Value = 3335.6667
$ws.Range("I99").Value = 2834.75
$ws.Range("J99").Value = 4337.5
$ws.Range("K99").Value = 2834.75
$ws.Range("L99").Value = 4337.5
$ws.Range("M99").Value = -1336.75
$ws.Range("N99").Value = -7333.5
$ws.Range("H107").Value = 1261.56
$ws.Range("I107").Value = 471.7857
$ws.Range("K107").Value = 471.7857
$ws.Range("M107").Value = 1448.2143
$ws.Range("H126").Value = 3335.6667
$ws.Range("I126").Value = 2834.75
$ws.Range("J126").Value = 4337.5
$ws.Range("K126").Value = 8504.25
$ws.Range("L126").Value = 13012.5
$ws.Range("M126").Value = -6034.25
$ws.Range("N126").Value = -17952.5
$ws.Range("I136").Value = 1859.6364
$ws.Range("J136").Value = 73104
$ws.Range("K136").Value = 5578.9092
$ws.Range("L136").Value = 219312
$ws.Range("M136").Value = -3028.9092
$ws.Range("N136").Value = -224412
$ws.Range("H141").Value = 11646.277
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 11646.277
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 11646.277
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -22006.277

# ---- CUL sheet ----
$ws = $wb.Worksheets("CUL")
$ws.Range("H24").Value = 659.8333
$ws.Range("I24").Value = 400
$ws.Range("J24").Value = 789.75
$ws.Range("K24").Value = 1200
$ws.Range("L24").Value = 2369.25
$ws.Range("M24").Value = -970
$ws.Range("N24").Value = -2829.25
$ws.Range("H69").Value = 1949.5
$ws.Range("I69").Value = 1500
$ws.Range("J69").Value = 2061.875
$ws.Range("K69").Value = 4500
$ws.Range("L69").Value = 6185.625
$ws.Range("M69").Value = -3689
$ws.Range("N69").Value = -7807.625
$ws.Range("H72").Value = 1949.5
$ws.Range("I72").Value = 1500
$ws.Range("J72").Value = 2061.875
$ws.Range("K72").Value = 13500
$ws.Range("L72").Value = 18556.875
$ws.Range("M72").Value = -9444
$ws.Range("N72").Value = -26668.875
$ws.Range("H109").Value = 3836
$ws.Range("I109").Value = 640.5714
$ws.Range("J109").Value = 5433.7144
$ws.Range("K109").Value = 1921.7142
$ws.Range("L109").Value = 16301.1432
$ws.Range("M109").Value = -881.7142000000001
$ws.Range("N109").Value = -18381.1432
$ws.Range("H117").Value = 1009.2
$ws.Range("I117").Value = 934.25
$ws.Range("J117").Value = 1036.4546
$ws.Range("K117").Value = 2802.75
$ws.Range("L117").Value = 3109.3638
$ws.Range("M117").Value = 639.25
$ws.Range("N117").Value = -9993.363799999999
$ws.Range("H131").Value = 691.72
$ws.Range("J131").Value = 719.47253
$ws.Range("L131").Value = 2158.41759
$ws.Range("N131").Value = -12238.41759

# ---- GSM sheet ----
$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 5220359
$ws.Range("I70").Value = 4629.5713
$ws.Range("J70").Value = 12522380
$ws.Range("K70").Value = 4629.5713
$ws.Range("L70").Value = 12522380
$ws.Range("M70").Value = -4359.5713
$ws.Range("N70").Value = -12522920
$ws.Range("H73").Value = 5220359
$ws.Range("I73").Value = 4629.5713
$ws.Range("J73").Value = 12522380
$ws.Range("K73").Value = 4629.5713
$ws.Range("L73").Value = 12522380
$ws.Range("M73").Value = -3693.5713
$ws.Range("N73").Value = -12524252
$ws.Range("H126").Value = 3004.9807
$ws.Range("I126").Value = 2492.6562
$ws.Range("K126").Value = 7477.9686
$ws.Range("M126").Value = -5007.9686
$ws.Range("H132").Value = 29725.3
$ws.Range("I132").Value = 4971.0586
$ws.Range("J132").Value = 169999.33
$ws.Range("K132").Value = 14913.1758
$ws.Range("L132").Value = 509997.99
$ws.Range("M132").Value = -12383.1758
$ws.Range("N132").Value = -515057.99

# ---- LTW sheet ----
$ws = $wb.Worksheets("LTW")
$ws.Range("H22").Value = 3049.3572
$ws.Range("I22").Value = 3365.9167
$ws.Range("J22").Value = 1150
$ws.Range("K22").Value = 3365.9167
$ws.Range("L22").Value = 1150
$ws.Range("M22").Value = -3070.9167
$ws.Range("N22").Value = -1740
$ws.Range("H27").Value = 3049.3572
$ws.Range("I27").Value = 3365.9167
$ws.Range("J27").Value = 1150
$ws.Range("K27").Value = 3365.9167
$ws.Range("L27").Value = 1150
$ws.Range("M27").Value = -3258.9167
$ws.Range("N27").Value = -1364
$ws.Range("H55").Value = 1040.4
$ws.Range("I55").Value = 1585
$ws.Range("J55").Value = 223.5
$ws.Range("K55").Value = 1585
$ws.Range("L55").Value = 223.5
$ws.Range("M55").Value = -1412
$ws.Range("N55").Value = -569.5

# ---- WVR sheet ----
$ws = $wb.Worksheets("WVR")
$ws.Range("H113").Value = 914.6316
$ws.Range("I113").Value = 961
$ws.Range("J113").Value = 80
$ws.Range("K113").Value = 2883
$ws.Range("L113").Value = 240
$ws.Range("M113").Value = -713
$ws.Range("N113").Value = -4580
$ws.Range("H139").Value = 51326.11
$ws.Range("J139").Value = 51326.11
$ws.Range("L139").Value = 51326.11
$ws.Range("N139").Value = -61606.11
